# data analysis.docx - apply the commit's edits via the Word object model.
$d = $word.ActiveDocument

# 1) Join the three runs ("...correlation ", "between", " the features.")
#    -- split apart only by spell/gram proofing marks -- back into a single
#    run/sentence. A scoped Find/Replace over the already-equal text lets
#    Word's own run-coalescing fold the (formatting-identical) runs back
#    into one run and drop the now-orphaned w:proofErr marks.
$d.Content.Find.Execute(
    " we can figure out the correlation between the features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " we can figure out the correlation between the features.", 2) | Out-Null

# 2) Retitle three of the four "BHK_OR_RK Vs TARGET PRICE BARPLOT" headings
#    (the first occurrence - the one right after the "ready to move" bullet,
#    which carries a lastRenderedPageBreak - is left alone).
#    Each heading paragraph holds the label in its own run, immediately
#    followed by a separate " Vs TARGET PRICE BARPLOT" run with the exact
#    same character formatting. A plain Range.Text replace on just the first
#    run would make Word's run-coalescer fuse the two back into one run, so
#    instead we briefly flip Bold off on the trailing run (forcing the two
#    runs to stay distinct across the edit), rewrite the label text, then
#    flip Bold back on -- leaving both runs exactly as formatted before.
$targets = @("BHK_NO.", "POSTED_BY", "SQUARE_FT")
$hit = 0
$targetIdx = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "BHK_OR_RK*") {
        $hit = $hit + 1
        if ($hit -gt 1) {
            $newLabel = $targets[$targetIdx]
            $targetIdx = $targetIdx + 1

            $start = $p.Range.Start
            $end = $p.Range.End

            $trailingRun = $d.Range($start + 9, $end)
            $trailingRun.Bold = 0

            $labelRun = $d.Range($start, $start + 9)
            $labelRun.Text = $newLabel

            $newTrailingStart = $start + $newLabel.Length
            $trailingRun2 = $d.Range($newTrailingStart, $p.Range.End)
            $trailingRun2.Bold = -1
        }
    }
}
